$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 19 (OpenRefine): mark as also NoSQL -> Model becomes Multi-Model
$ws.Range("D19").Value = $true
$ws.Range("E19").Value = "Multi-Model"

# Row 156 (nutch): mark as also NoSQL -> Model becomes Multi-Model
$ws.Range("D156").Value = $true
$ws.Range("E156").Value = "Multi-Model"

# Insert a new project row before the old row 224 ("vert.x"), shifting
# rows 224-233 down to 225-234 (zookeeper ends up on row 234).
$ws.Rows(224).Insert()

# Match the bordered/bold style used by column A on the other data rows.
$ws.Range("A223").Copy($ws.Range("A224"))

$ws.Range("A224").Value = 346
$ws.Range("B224").Value = "validator"
$ws.Range("C224").Value = $false
$ws.Range("D224").Value = $true
$ws.Range("E224").Value = "NoSQL"
